# Append a new "2025-03-12" price row (row 11) to every price sheet in the
# workbook, mirroring the last recorded price (row 10) for each series,
# except USD_CNY which already had its own fresh value.
#
# Values/dates are entered as plain text (matching the existing sheet
# convention where every other cell in these columns is stored as text,
# not as a real number/date), so we force the cell to Text format before
# assigning the value and then restore the "Normal" style so no stray
# number-format is left applied to the cell.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-12"

$sheetValues = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.295"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,221"
    "Silver Busbar front-side"   = "7,816"
    "Silver finger front-side"   = "7,866"
    "USD_CNY"                    = "7.2787"
}

foreach ($sheetName in $sheetValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $price = $sheetValues[$sheetName]

    $dateCell = $ws.Cells.Item(11, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate
    $dateCell.Style = "Normal"

    $priceCell = $ws.Cells.Item(11, 2)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"
}
